$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foods")

# Update existing rows 2 and 3, and add new rows 4..21
# Columns: A=idFood, B=Name, C=AlternativeName (kept blank), D=Description
$data = @(
    @("F01", "Apples",        "Fresh apples"),
    @("F02", "Pears",         "Fresh pears"),
    @("F03", "Strawberries",  "Fresh strawberries"),
    @("F04", "Grapes",        "Fresh grapes"),
    @("F05", "Tomatoes",      "Fresh tomatoes"),
    @("F06", "Lettuce",       "Leafy lettuce"),
    @("F07", "Spinach",       "Fresh spinach"),
    @("F08", "Potatoes",      "Raw potatoes"),
    @("F09", "Wheat bread",   "Wheat bread"),
    @("F10", "Rice",          "Dry rice"),
    @("F11", "Oats",          "Rolled oats"),
    @("F12", "Milk",          "Whole milk"),
    @("F13", "Yogurt",        "Plain yogurt"),
    @("F14", "Chicken",       "Chicken meat"),
    @("F15", "Pork",          "Pork meat"),
    @("F16", "Salmon",        "Salmon fillet"),
    @("F17", "Eggs",          "Chicken eggs"),
    @("F18", "Olive oil",     "Extra virgin olive oil"),
    @("F19", "Oranges",       "Fresh oranges"),
    @("F20", "Carrots",       "Fresh carrots")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 4).Value = $rowData[2]
}

# Column C (AlternativeName) stays blank for every row, same as rows 2-3 in the
# original sheet. Propagate the existing blank cell down to the new rows so the
# cells remain present (instead of disappearing when assigned "").
$ws.Range("C2").Copy($ws.Range("C4:C21"))
